# "perbaikan italic pada bahasa asing" - italicize the foreign (English) word
# "website" / "WEBSITE" everywhere it appears as a standalone loanword in the
# running Indonesian text (one occurrence was already italic and is left
# untouched).
$d = $word.ActiveDocument

function Italicize-Word {
    param(
        [string]$OuterPhrase,
        [string]$InnerWord
    )

    $full = $d.Content
    $ok = $full.Find.Execute($OuterPhrase, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "outer phrase not found: $OuterPhrase"
    }

    $local = $d.Range($full.Start, $full.End)
    $ok2 = $local.Find.Execute($InnerWord, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok2) {
        throw "inner word not found: $InnerWord (within $OuterPhrase)"
    }

    $target = $d.Range($local.Start, $local.End)
    $target.Italic = 1
}

# 1. "...berkembang adalah website, yang telah menjadi wadah..."
Italicize-Word "adalah website, yang telah menjadi" "website"

# 2. "...media utama pembangun website Bahasa pemrograman javascript merupakan..."
Italicize-Word "pembangun website Bahasa pemrograman" "website"

# 3. "...untuk membangun tampilan website secara interaktif..."
Italicize-Word "membangun tampilan website secara interaktif" "website"

# 4. "...Pengembangan front-end dan back-end website ini menjadi fokus..."
Italicize-Word "back-end website ini menjadi fokus" "website"

# 5. "...serta mengelola website secara profesional sesuai..."
Italicize-Word "mengelola website secara profesional" "website"

# 6. Title: "...FRONT-END DAN BACK-END WEBSITE NUSANTARAKU..."
Italicize-Word "BACK-END WEBSITE NUSANTARAKU" "WEBSITE"

# 7. "...mengenalkan budaya di seluruh Indonesia pada website NusantaraKu"
Italicize-Word "pada website NusantaraKu" "website"

# 8. "...berbasis website menggunakan bahasa pemograman PHP..."
Italicize-Word "berbasis website menggunakan bahasa" "website"

# 9. "...sistem,konsep sistem informasi,website, database..."
Italicize-Word "informasi,website, database" "website"

# 10. "...pengujian sistem informasi website yang telah dibangun..."
Italicize-Word "pengujian sistem informasi website yang telah dibangun" "website"

Write-Output "done"
